$wb = $excel.ActiveWorkbook

# --- Sheet "contest": update summary stats for Weekly Contest 321 (row 2) ---
$contest = $wb.Worksheets.Item("contest")
$contest.Cells.Item(2, 2).Value = 0       # B2 2q_sum
$contest.Cells.Item(2, 3).Value = 0       # C2 2q_pop
$contest.Cells.Item(2, 6).Value = 75790    # F2 4q_sum
$contest.Cells.Item(2, 7).Value = 25       # G2 4q_pop

# --- Sheet "user": rewrite leaderboard with new columns and rows ---
$ws = $wb.Worksheets.Item("user")

# New header cells for columns I..M
$ws.Cells.Item(1, 9).Value = "views"
$ws.Cells.Item(1, 10).Value = "solution"
$ws.Cells.Item(1, 11).Value = "discuss"
$ws.Cells.Item(1, 12).Value = "reputation"
$ws.Cells.Item(1, 13).Value = "reput_level"

# Match header style (bold, centered, bordered) used by existing header cells
$ws.Range("H1").Copy()
$ws.Range("I1:M1").PasteSpecial(-4122)

# Row 2: JOHNKRAM
$ws.Cells.Item(2, 1).Value = "JOHNKRAM"
$ws.Cells.Item(2, 2).Value = 3584
$ws.Cells.Item(2, 3).Value = "China"
$ws.Cells.Item(2, 4).Value = ""
$ws.Cells.Item(2, 5).Value = ""
$ws.Cells.Item(2, 6).Value = "清华大学"
$ws.Cells.Item(2, 7).Value = "C++"
$ws.Cells.Item(2, 8).Value = 1
$ws.Cells.Item(2, 9).Value = 19400
$ws.Cells.Item(2, 10).Value = ""
$ws.Cells.Item(2, 11).Value = ""
$ws.Cells.Item(2, 12).Value = ""
$ws.Cells.Item(2, 13).Value = "L1"

# Row 3: qeetcode
$ws.Cells.Item(3, 1).Value = "qeetcode"
$ws.Cells.Item(3, 2).Value = 3257
$ws.Cells.Item(3, 3).Value = "United States"
$ws.Cells.Item(3, 4).Value = ""
$ws.Cells.Item(3, 5).Value = ""
$ws.Cells.Item(3, 6).Value = "University of California--Berkeley"
$ws.Cells.Item(3, 7).Value = "C++"
$ws.Cells.Item(3, 8).Value = 1
$ws.Cells.Item(3, 9).Value = 119000
$ws.Cells.Item(3, 10).Value = 972
$ws.Cells.Item(3, 11).Value = 0
$ws.Cells.Item(3, 12).Value = 578
$ws.Cells.Item(3, 13).Value = ""

# Row 4: bucketpotato
$ws.Cells.Item(4, 1).Value = "bucketpotato"
$ws.Cells.Item(4, 2).Value = 3169
$ws.Cells.Item(4, 3).Value = "United States"
$ws.Cells.Item(4, 4).Value = ""
$ws.Cells.Item(4, 5).Value = ""
$ws.Cells.Item(4, 6).Value = ""
$ws.Cells.Item(4, 7).Value = "C++"
$ws.Cells.Item(4, 8).Value = 1
$ws.Cells.Item(4, 9).Value = 0
$ws.Cells.Item(4, 10).Value = 0
$ws.Cells.Item(4, 11).Value = 0
$ws.Cells.Item(4, 12).Value = 0
$ws.Cells.Item(4, 13).Value = ""

# Row 5: moransky
$ws.Cells.Item(5, 1).Value = "moransky"
$ws.Cells.Item(5, 2).Value = 3368
$ws.Cells.Item(5, 3).Value = "China"
$ws.Cells.Item(5, 4).Value = ""
$ws.Cells.Item(5, 5).Value = ""
$ws.Cells.Item(5, 6).Value = ""
$ws.Cells.Item(5, 7).Value = "C++"
$ws.Cells.Item(5, 8).Value = 1
$ws.Cells.Item(5, 9).Value = 0
$ws.Cells.Item(5, 10).Value = ""
$ws.Cells.Item(5, 11).Value = ""
$ws.Cells.Item(5, 12).Value = ""
$ws.Cells.Item(5, 13).Value = "暂无"

# Row 6: jinmingli
$ws.Cells.Item(6, 1).Value = "jinmingli"
$ws.Cells.Item(6, 2).Value = 2997
$ws.Cells.Item(6, 3).Value = "China"
$ws.Cells.Item(6, 4).Value = "高德地图"
$ws.Cells.Item(6, 5).Value = "算法专家"
$ws.Cells.Item(6, 6).Value = "清华大学"
$ws.Cells.Item(6, 7).Value = "C++"
$ws.Cells.Item(6, 8).Value = 1
$ws.Cells.Item(6, 9).Value = 0
$ws.Cells.Item(6, 10).Value = ""
$ws.Cells.Item(6, 11).Value = ""
$ws.Cells.Item(6, 12).Value = ""
$ws.Cells.Item(6, 13).Value = "暂无"

# Row 7: galencolin
$ws.Cells.Item(7, 1).Value = "galencolin"
$ws.Cells.Item(7, 2).Value = 2817
$ws.Cells.Item(7, 3).Value = "Unknown"
$ws.Cells.Item(7, 4).Value = ""
$ws.Cells.Item(7, 5).Value = ""
$ws.Cells.Item(7, 6).Value = ""
$ws.Cells.Item(7, 7).Value = "C++"
$ws.Cells.Item(7, 8).Value = 1
$ws.Cells.Item(7, 9).Value = 0
$ws.Cells.Item(7, 10).Value = 0
$ws.Cells.Item(7, 11).Value = 0
$ws.Cells.Item(7, 12).Value = 0
$ws.Cells.Item(7, 13).Value = ""

# Row 8: lucifer1006
$ws.Cells.Item(8, 1).Value = "lucifer1006"
$ws.Cells.Item(8, 2).Value = 3097
$ws.Cells.Item(8, 3).Value = "China"
$ws.Cells.Item(8, 4).Value = "Viktor Chondria University"
$ws.Cells.Item(8, 5).Value = "Researcher"
$ws.Cells.Item(8, 6).Value = "北京大学"
$ws.Cells.Item(8, 7).Value = "C++"
$ws.Cells.Item(8, 8).Value = 1
$ws.Cells.Item(8, 9).Value = 0
$ws.Cells.Item(8, 10).Value = ""
$ws.Cells.Item(8, 11).Value = ""
$ws.Cells.Item(8, 12).Value = ""
$ws.Cells.Item(8, 13).Value = "暂无"

# Row 9: nyu_ldf
$ws.Cells.Item(9, 1).Value = "nyu_ldf"
$ws.Cells.Item(9, 2).Value = 3494
$ws.Cells.Item(9, 3).Value = "Unknown"
$ws.Cells.Item(9, 4).Value = ""
$ws.Cells.Item(9, 5).Value = ""
$ws.Cells.Item(9, 6).Value = ""
$ws.Cells.Item(9, 7).Value = "Python"
$ws.Cells.Item(9, 8).Value = 1
$ws.Cells.Item(9, 9).Value = 0
$ws.Cells.Item(9, 10).Value = 0
$ws.Cells.Item(9, 11).Value = 0
$ws.Cells.Item(9, 12).Value = 0
$ws.Cells.Item(9, 13).Value = ""

# Row 10: liouzhou_101
$ws.Cells.Item(10, 1).Value = "liouzhou_101"
$ws.Cells.Item(10, 2).Value = 3204
$ws.Cells.Item(10, 3).Value = "China"
$ws.Cells.Item(10, 4).Value = ""
$ws.Cells.Item(10, 5).Value = ""
$ws.Cells.Item(10, 6).Value = "清华大学"
$ws.Cells.Item(10, 7).Value = "C++"
$ws.Cells.Item(10, 8).Value = 1
$ws.Cells.Item(10, 9).Value = 0
$ws.Cells.Item(10, 10).Value = ""
$ws.Cells.Item(10, 11).Value = ""
$ws.Cells.Item(10, 12).Value = ""
$ws.Cells.Item(10, 13).Value = "暂无"

# Row 11: c8kbf
$ws.Cells.Item(11, 1).Value = "c8kbf"
$ws.Cells.Item(11, 2).Value = 2923
$ws.Cells.Item(11, 3).Value = "Canada"
$ws.Cells.Item(11, 4).Value = ""
$ws.Cells.Item(11, 5).Value = ""
$ws.Cells.Item(11, 6).Value = ""
$ws.Cells.Item(11, 7).Value = "C++"
$ws.Cells.Item(11, 8).Value = 1
$ws.Cells.Item(11, 9).Value = 0
$ws.Cells.Item(11, 10).Value = 2
$ws.Cells.Item(11, 11).Value = 0
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = ""

# Row 12: raincoat911
$ws.Cells.Item(12, 1).Value = "raincoat911"
$ws.Cells.Item(12, 2).Value = 2901
$ws.Cells.Item(12, 3).Value = "United States"
$ws.Cells.Item(12, 4).Value = ""
$ws.Cells.Item(12, 5).Value = ""
$ws.Cells.Item(12, 6).Value = ""
$ws.Cells.Item(12, 7).Value = "C++"
$ws.Cells.Item(12, 8).Value = 1
$ws.Cells.Item(12, 9).Value = 405
$ws.Cells.Item(12, 10).Value = 1
$ws.Cells.Item(12, 11).Value = 0
$ws.Cells.Item(12, 12).Value = 5
$ws.Cells.Item(12, 13).Value = ""

# Row 13: delphih
$ws.Cells.Item(13, 1).Value = "delphih"
$ws.Cells.Item(13, 2).Value = 2768
$ws.Cells.Item(13, 3).Value = "United States"
$ws.Cells.Item(13, 4).Value = ""
$ws.Cells.Item(13, 5).Value = ""
$ws.Cells.Item(13, 6).Value = "Georgia Institute of Technology"
$ws.Cells.Item(13, 7).Value = "Python3"
$ws.Cells.Item(13, 8).Value = 1
$ws.Cells.Item(13, 9).Value = 16400
$ws.Cells.Item(13, 10).Value = 29
$ws.Cells.Item(13, 11).Value = 0
$ws.Cells.Item(13, 12).Value = 219
$ws.Cells.Item(13, 13).Value = ""

# Row 14: arignote
$ws.Cells.Item(14, 1).Value = "arignote"
$ws.Cells.Item(14, 2).Value = 3408
$ws.Cells.Item(14, 3).Value = "China"
$ws.Cells.Item(14, 4).Value = ""
$ws.Cells.Item(14, 5).Value = ""
$ws.Cells.Item(14, 6).Value = "海外高校"
$ws.Cells.Item(14, 7).Value = "Java"
$ws.Cells.Item(14, 8).Value = 1
$ws.Cells.Item(14, 9).Value = 0
$ws.Cells.Item(14, 10).Value = ""
$ws.Cells.Item(14, 11).Value = ""
$ws.Cells.Item(14, 12).Value = ""
$ws.Cells.Item(14, 13).Value = "暂无"

# Row 15: lympanda
$ws.Cells.Item(15, 1).Value = "lympanda"
$ws.Cells.Item(15, 2).Value = 3139
$ws.Cells.Item(15, 3).Value = "Unknown"
$ws.Cells.Item(15, 4).Value = ""
$ws.Cells.Item(15, 5).Value = ""
$ws.Cells.Item(15, 6).Value = ""
$ws.Cells.Item(15, 7).Value = "C++"
$ws.Cells.Item(15, 8).Value = 1
$ws.Cells.Item(15, 9).Value = 220
$ws.Cells.Item(15, 10).Value = 3
$ws.Cells.Item(15, 11).Value = 0
$ws.Cells.Item(15, 12).Value = 0
$ws.Cells.Item(15, 13).Value = ""

# Row 16: cpp_template
$ws.Cells.Item(16, 1).Value = "cpp_template"
$ws.Cells.Item(16, 2).Value = 3104
$ws.Cells.Item(16, 3).Value = "China"
$ws.Cells.Item(16, 4).Value = ""
$ws.Cells.Item(16, 5).Value = ""
$ws.Cells.Item(16, 6).Value = ""
$ws.Cells.Item(16, 7).Value = "C++"
$ws.Cells.Item(16, 8).Value = 1
$ws.Cells.Item(16, 9).Value = 0
$ws.Cells.Item(16, 10).Value = ""
$ws.Cells.Item(16, 11).Value = ""
$ws.Cells.Item(16, 12).Value = ""
$ws.Cells.Item(16, 13).Value = "暂无"

# Row 17: Yawn_Sean
$ws.Cells.Item(17, 1).Value = "Yawn_Sean"
$ws.Cells.Item(17, 2).Value = 3118
$ws.Cells.Item(17, 3).Value = "China"
$ws.Cells.Item(17, 4).Value = ""
$ws.Cells.Item(17, 5).Value = ""
$ws.Cells.Item(17, 6).Value = "北京大学"
$ws.Cells.Item(17, 7).Value = "Python3"
$ws.Cells.Item(17, 8).Value = 1
$ws.Cells.Item(17, 9).Value = 0
$ws.Cells.Item(17, 10).Value = ""
$ws.Cells.Item(17, 11).Value = ""
$ws.Cells.Item(17, 12).Value = ""
$ws.Cells.Item(17, 13).Value = "L1"

# Row 18: Tlatoani
$ws.Cells.Item(18, 1).Value = "Tlatoani"
$ws.Cells.Item(18, 2).Value = 3225
$ws.Cells.Item(18, 3).Value = "Mexico"
$ws.Cells.Item(18, 4).Value = ""
$ws.Cells.Item(18, 5).Value = ""
$ws.Cells.Item(18, 6).Value = ""
$ws.Cells.Item(18, 7).Value = "Kotlin"
$ws.Cells.Item(18, 8).Value = 1
$ws.Cells.Item(18, 9).Value = 0
$ws.Cells.Item(18, 10).Value = 0
$ws.Cells.Item(18, 11).Value = 0
$ws.Cells.Item(18, 12).Value = 0
$ws.Cells.Item(18, 13).Value = ""

# Row 19: user3754Ay
$ws.Cells.Item(19, 1).Value = "user3754Ay"
$ws.Cells.Item(19, 2).Value = 2027
$ws.Cells.Item(19, 3).Value = "Unknown"
$ws.Cells.Item(19, 4).Value = ""
$ws.Cells.Item(19, 5).Value = ""
$ws.Cells.Item(19, 6).Value = "Massachusetts Institute of Technology"
$ws.Cells.Item(19, 7).Value = "C++"
$ws.Cells.Item(19, 8).Value = 1
$ws.Cells.Item(19, 9).Value = 0
$ws.Cells.Item(19, 10).Value = 0
$ws.Cells.Item(19, 11).Value = 0
$ws.Cells.Item(19, 12).Value = 0
$ws.Cells.Item(19, 13).Value = ""

# Row 20: bitetheD4T
$ws.Cells.Item(20, 1).Value = "bitetheD4T"
$ws.Cells.Item(20, 2).Value = 2756
$ws.Cells.Item(20, 3).Value = "China"
$ws.Cells.Item(20, 4).Value = ""
$ws.Cells.Item(20, 5).Value = ""
$ws.Cells.Item(20, 6).Value = "浙大宁波理工学院"
$ws.Cells.Item(20, 7).Value = "C++"
$ws.Cells.Item(20, 8).Value = 1
$ws.Cells.Item(20, 9).Value = 563
$ws.Cells.Item(20, 10).Value = ""
$ws.Cells.Item(20, 11).Value = ""
$ws.Cells.Item(20, 12).Value = ""
$ws.Cells.Item(20, 13).Value = "L1"

# Row 21: jianghd1996
$ws.Cells.Item(21, 1).Value = "jianghd1996"
$ws.Cells.Item(21, 2).Value = 2712
$ws.Cells.Item(21, 3).Value = "China"
$ws.Cells.Item(21, 4).Value = ""
$ws.Cells.Item(21, 5).Value = ""
$ws.Cells.Item(21, 6).Value = "北京大学"
$ws.Cells.Item(21, 7).Value = "Python"
$ws.Cells.Item(21, 8).Value = 1
$ws.Cells.Item(21, 9).Value = 0
$ws.Cells.Item(21, 10).Value = ""
$ws.Cells.Item(21, 11).Value = ""
$ws.Cells.Item(21, 12).Value = ""
$ws.Cells.Item(21, 13).Value = "暂无"

# Row 22: lxhgww
$ws.Cells.Item(22, 1).Value = "lxhgww"
$ws.Cells.Item(22, 2).Value = 2772
$ws.Cells.Item(22, 3).Value = "China"
$ws.Cells.Item(22, 4).Value = ""
$ws.Cells.Item(22, 5).Value = ""
$ws.Cells.Item(22, 6).Value = "清华大学"
$ws.Cells.Item(22, 7).Value = "C++"
$ws.Cells.Item(22, 8).Value = 1
$ws.Cells.Item(22, 9).Value = 0
$ws.Cells.Item(22, 10).Value = ""
$ws.Cells.Item(22, 11).Value = ""
$ws.Cells.Item(22, 12).Value = ""
$ws.Cells.Item(22, 13).Value = "暂无"

# Row 23: sammochen
$ws.Cells.Item(23, 1).Value = "sammochen"
$ws.Cells.Item(23, 2).Value = 3049
$ws.Cells.Item(23, 3).Value = "New Zealand"
$ws.Cells.Item(23, 4).Value = ""
$ws.Cells.Item(23, 5).Value = ""
$ws.Cells.Item(23, 6).Value = "University of Auckland"
$ws.Cells.Item(23, 7).Value = "C++"
$ws.Cells.Item(23, 8).Value = 1
$ws.Cells.Item(23, 9).Value = 0
$ws.Cells.Item(23, 10).Value = 1
$ws.Cells.Item(23, 11).Value = 0
$ws.Cells.Item(23, 12).Value = 0
$ws.Cells.Item(23, 13).Value = ""

# Row 24: w285714
$ws.Cells.Item(24, 1).Value = "w285714"
$ws.Cells.Item(24, 2).Value = 2938
$ws.Cells.Item(24, 3).Value = "China"
$ws.Cells.Item(24, 4).Value = ""
$ws.Cells.Item(24, 5).Value = ""
$ws.Cells.Item(24, 6).Value = "北京大学"
$ws.Cells.Item(24, 7).Value = "C++"
$ws.Cells.Item(24, 8).Value = 1
$ws.Cells.Item(24, 9).Value = 0
$ws.Cells.Item(24, 10).Value = ""
$ws.Cells.Item(24, 11).Value = ""
$ws.Cells.Item(24, 12).Value = ""
$ws.Cells.Item(24, 13).Value = "暂无"

# Row 25: uwi
$ws.Cells.Item(25, 1).Value = "uwi"
$ws.Cells.Item(25, 2).Value = 3463
$ws.Cells.Item(25, 3).Value = "Japan"
$ws.Cells.Item(25, 4).Value = ""
$ws.Cells.Item(25, 5).Value = ""
$ws.Cells.Item(25, 6).Value = ""
$ws.Cells.Item(25, 7).Value = "Java"
$ws.Cells.Item(25, 8).Value = 1
$ws.Cells.Item(25, 9).Value = 787
$ws.Cells.Item(25, 10).Value = 1
$ws.Cells.Item(25, 11).Value = 0
$ws.Cells.Item(25, 12).Value = 21
$ws.Cells.Item(25, 13).Value = ""

# Row 26: ray_striker
$ws.Cells.Item(26, 1).Value = "ray_striker"
$ws.Cells.Item(26, 2).Value = 2500
$ws.Cells.Item(26, 3).Value = "India"
$ws.Cells.Item(26, 4).Value = ""
$ws.Cells.Item(26, 5).Value = ""
$ws.Cells.Item(26, 6).Value = "National Institute of Technology, Silchar"
$ws.Cells.Item(26, 7).Value = "C++"
$ws.Cells.Item(26, 8).Value = 1
$ws.Cells.Item(26, 9).Value = 5100
$ws.Cells.Item(26, 10).Value = 4
$ws.Cells.Item(26, 11).Value = 1
$ws.Cells.Item(26, 12).Value = 9
$ws.Cells.Item(26, 13).Value = ""

# New rows 19..26 need column A to match the bold/bordered style used elsewhere in column A
$ws.Range("A2").Copy()
$ws.Range("A19:A26").PasteSpecial(-4122)
$excel.CutCopyMode = $false

